# =====================================================================
# Edit script: restructure PlayerPerformance workbook
#  1. Insert new "Player Info" sheet at the front
#  2. Rename ODI Batting's MATCH_CARD_LINK column to MATCH_CODE and
#     replace the URL values with the bare numeric match code
#  3. Rename ODI Bowling's MATCH_CARD_LINK column to MATCH_CODE and
#     replace the URL values with the bare numeric match code
#  4. Append new "ODI Batting Extra" sheet at the end
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Player Info" sheet (inserted before the first existing sheet)
# ---------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Count; $c++) {
    $playerInfo.Cells.Item(1, $c).Value2 = $piHeaders[$c - 1]
}
$piHdrRange = $playerInfo.Range($playerInfo.Cells.Item(1, 1), $playerInfo.Cells.Item(1, $piHeaders.Count))
$piHdrRange.Font.Bold = $true
$piHdrRange.HorizontalAlignment = -4108
$piHdrRange.VerticalAlignment = -4160
$piHdrRange.Borders.LineStyle = 1

$piDataRange = $playerInfo.Range($playerInfo.Cells.Item(2, 1), $playerInfo.Cells.Item(2, $piHeaders.Count))
$piDataRange.NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value2 = "3665"
$playerInfo.Cells.Item(2, 2).Value2 = "Mohammad Nabi"
$playerInfo.Cells.Item(2, 3).Value2 = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value2 = "Right Arm Off Break"

# ---------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$battingLastRow = $batting.UsedRange.Rows.Count
$batting.Cells.Item(1, 4).Value2 = "MATCH_CODE"

$battingCodeRange = $batting.Range($batting.Cells.Item(2, 4), $batting.Cells.Item($battingLastRow, 4))
$battingCodeRange.NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null) {
        $code = $v -replace ".*MatchCode=(\d+).*", '$1'
        $cell.Value2 = $code
    }
}

# ---------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlingLastRow = $bowling.UsedRange.Rows.Count
$bowling.Cells.Item(1, 2).Value2 = "MATCH_CODE"

$bowlingCodeRange = $bowling.Range($bowling.Cells.Item(2, 2), $bowling.Cells.Item($bowlingLastRow, 2))
$bowlingCodeRange.NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -ne $null) {
        $code = $v -replace ".*MatchCode=(\d+).*", '$1'
        $cell.Value2 = $code
    }
}

# ---------------------------------------------------------------
# 4. "ODI Batting Extra" sheet (appended after the last existing sheet)
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$beHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $beHeaders.Count; $c++) {
    $battingExtra.Cells.Item(1, $c).Value2 = $beHeaders[$c - 1]
}
$beHdrRange = $battingExtra.Range($battingExtra.Cells.Item(1, 1), $battingExtra.Cells.Item(1, $beHeaders.Count))
$beHdrRange.Font.Bold = $true
$beHdrRange.HorizontalAlignment = -4108
$beHdrRange.VerticalAlignment = -4160
$beHdrRange.Borders.LineStyle = 1

$beData = @(
    @("4326", $null, $null, $null, $null, "NO"),
    @("4332", $null, $null, $null, $null, "NO"),
    @("4335", 5, "0", "0", $null, "NO"),
    @("4340", 6, "0", "0", "7.05%", "NO"),
    @("4348", 6, "0", "0", "0.69%", "NO"),
    @("4377", 7, "0", "0", "0.52%", "NO"),
    @("4378", 7, "5", "0", "16.00%", "NO"),
    @("4379", 7, "3", "1", "20.08%", "NO"),
    @("4444", 5, "0", "0", "1.05%", "NO"),
    @("4446", 7, $null, $null, $null, "NO"),
    @("4448", 6, "3", "1", "12.03%", "NO"),
    @("4537", 6, "2", "0", "9.30%", "NO"),
    @("4538", $null, $null, $null, $null, "NO"),
    @("4539", 6, $null, $null, $null, "NO"),
    @("4582", 5, "1", "0", "3.62%", "NO"),
    @("4585", 6, $null, $null, $null, "NO"),
    @("4588", 6, "2", "1", "24.82%", "NO"),
    @("4671", $null, $null, $null, $null, "NO"),
    @("4674", 7, "2", "2", "17.98%", "NO"),
    @("4675", $null, $null, $null, $null, "NO")
)

$beLastRow = $beData.Count + 1
$beTextCols = @(1, 3, 4, 5, 6)
foreach ($col in $beTextCols) {
    $rng = $battingExtra.Range($battingExtra.Cells.Item(2, $col), $battingExtra.Cells.Item($beLastRow, $col))
    $rng.NumberFormat = "@"
}

$r = 2
foreach ($row in $beData) {
    for ($c = 1; $c -le $row.Count; $c++) {
        $val = $row[$c - 1]
        if ($val -ne $null) {
            $battingExtra.Cells.Item($r, $c).Value2 = $val
        }
    }
    $r++
}

$playerInfo.Activate()

Write-Host "Edit complete"
